# Change cell B11 (on the active/Rules sheet) from the shared text "R40" to
# the literal text "1". Assigning a numeric-looking string straight to
# Range.Value (e.g. "1") gets auto-coerced to a number by Excel, and forcing
# text via NumberFormat="@"/quote-prefix on B11 itself would also mutate its
# existing cell style (fillId/borderId) which must stay untouched (s="23").
#
# Workaround: stage the text value "1" in an unused scratch cell (A1, outside
# the sheet's B3:E11 used range) that's explicitly formatted as Text, then
# copy only the *value* (PasteSpecial xlPasteValues = -4163) into B11 so its
# existing number format/style is left completely alone. Finally wipe the
# scratch cell so no trace of it remains in the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("A1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"
$scratch.Copy()

$ws.Range("B11").PasteSpecial(-4163)

$scratch.Clear()
